$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).NumberFormat = "@"
$ws.Cells.Item(2, 4).Value = "48.803.34"
$ws.Cells.Item(2, 4).ClearFormats()
$ws.Range("E2").Value = "  -1.72%  "
$ws.Cells.Item(3, 4).NumberFormat = "@"
$ws.Cells.Item(3, 4).Value = "2.637.59"
$ws.Cells.Item(3, 4).ClearFormats()
$ws.Range("E3").Value = "  +3.20%  "
$ws.Range("E4").Value = "  +0.22%  "
$ws.Cells.Item(5, 4).NumberFormat = "@"
$ws.Cells.Item(5, 4).Value = "110.51"
$ws.Cells.Item(5, 4).ClearFormats()
$ws.Range("E5").Value = "  +1.88%  "
$ws.Cells.Item(6, 4).NumberFormat = "@"
$ws.Cells.Item(6, 4).Value = "322.67"
$ws.Cells.Item(6, 4).ClearFormats()
$ws.Range("E6").Value = "  +0.30%  "
$ws.Range("E7").Value = "  -1.33%  "
$ws.Range("E8").Value = "  +0.10%  "
$ws.Range("E9").Value = "  -2.42%  "
$ws.Cells.Item(10, 4).NumberFormat = "@"
$ws.Cells.Item(10, 4).Value = "39.51"
$ws.Cells.Item(10, 4).ClearFormats()
$ws.Range("E10").Value = "  -1.31%  "
$ws.Cells.Item(11, 4).NumberFormat = "@"
$ws.Cells.Item(11, 4).Value = "19.88"
$ws.Cells.Item(11, 4).ClearFormats()
$ws.Range("E11").Value = "  -2.17%  "
$ws.Cells.Item(12, 4).NumberFormat = "@"
$ws.Cells.Item(12, 4).Value = "0.0810"
$ws.Cells.Item(12, 4).ClearFormats()
$ws.Range("E12").Value = "  -0.43%  "
$ws.Range("E13").Value = "  -0.47%  "
$ws.Range("E14").Value = "  -0.45%  "
$ws.Cells.Item(15, 4).NumberFormat = "@"
$ws.Cells.Item(15, 4).Value = "3.045.54"
$ws.Cells.Item(15, 4).ClearFormats()
$ws.Range("E15").Value = "  +3.39%  "
$ws.Cells.Item(16, 4).NumberFormat = "@"
$ws.Cells.Item(16, 4).Value = "2.632.00"
$ws.Cells.Item(16, 4).ClearFormats()
$ws.Range("E16").Value = "  +2.74%  "
$ws.Cells.Item(17, 4).NumberFormat = "@"
$ws.Cells.Item(17, 4).Value = "0.864"
$ws.Cells.Item(17, 4).ClearFormats()
$ws.Range("E17").Value = "  +0.74%  "
$ws.Cells.Item(18, 4).NumberFormat = "@"
$ws.Cells.Item(18, 4).Value = "48.803.67"
$ws.Cells.Item(18, 4).ClearFormats()
$ws.Range("E18").Value = "  -1.31%  "
$ws.Cells.Item(19, 4).NumberFormat = "@"
$ws.Cells.Item(19, 4).Value = "12.85"
$ws.Cells.Item(19, 4).ClearFormats()
$ws.Range("E19").Value = "  -2.39%  "
$ws.Range("E20").Value = "  -0.40%  "
$ws.Range("E21").Value = "  -2.03%  "
$ws.Cells.Item(22, 4).NumberFormat = "@"
$ws.Cells.Item(22, 4).Value = "0.0₃0942"
$ws.Cells.Item(22, 4).ClearFormats()
$ws.Range("E22").Value = "  -0.15%  "
$ws.Cells.Item(23, 4).NumberFormat = "@"
$ws.Cells.Item(23, 4).Value = "270.11"
$ws.Cells.Item(23, 4).ClearFormats()
$ws.Range("E23").Value = "  -4.87%  "
$ws.Cells.Item(24, 4).NumberFormat = "@"
$ws.Cells.Item(24, 4).Value = "69.18"
$ws.Cells.Item(24, 4).ClearFormats()
$ws.Range("E24").Value = "  -3.63%  "
$ws.Cells.Item(25, 4).NumberFormat = "@"
$ws.Cells.Item(25, 4).Value = "2.52"
$ws.Cells.Item(25, 4).ClearFormats()
$ws.Range("E25").Value = "  +0.25%  "
$ws.Cells.Item(26, 4).NumberFormat = "@"
$ws.Cells.Item(26, 4).Value = "26.05"
$ws.Cells.Item(26, 4).ClearFormats()
$ws.Range("E26").Value = "  -1.27%  "
$ws.Range("E27").Value = "  +0.06%  "
$ws.Cells.Item(28, 4).NumberFormat = "@"
$ws.Cells.Item(28, 4).Value = "10.11"
$ws.Cells.Item(28, 4).ClearFormats()
$ws.Range("E28").Value = "  +3.13%  "
$ws.Range("E29").Value = "  +0.14%  "
$ws.Cells.Item(30, 4).NumberFormat = "@"
$ws.Cells.Item(30, 4).Value = "34.90"
$ws.Cells.Item(30, 4).ClearFormats()
$ws.Range("E30").Value = "  -1.32%  "
$ws.Range("E31").Value = "  -5.39%  "
$ws.Cells.Item(32, 4).NumberFormat = "@"
$ws.Cells.Item(32, 4).Value = "49.42"
$ws.Cells.Item(32, 4).ClearFormats()
$ws.Range("E32").Value = "  -0.38%  "
$ws.Range("E33").Value = "  +1.06%  "
$ws.Cells.Item(34, 4).NumberFormat = "@"
$ws.Cells.Item(34, 4).Value = "19.24"
$ws.Cells.Item(34, 4).ClearFormats()
$ws.Range("E34").Value = "  -2.22%  "
$ws.Range("E35").Value = "  -0.06%  "
$ws.Cells.Item(36, 4).NumberFormat = "@"
$ws.Cells.Item(36, 4).Value = "0.0796"
$ws.Cells.Item(36, 4).ClearFormats()
$ws.Range("E36").Value = "  +1.69%  "
$ws.Range("E37").Value = "  +5.50%  "
$ws.Range("E38").Value = "  +1.30%  "
$ws.Cells.Item(39, 4).NumberFormat = "@"
$ws.Cells.Item(39, 4).Value = "3.15"
$ws.Cells.Item(39, 4).ClearFormats()
$ws.Range("E39").Value = "  +6.24%  "
$ws.Cells.Item(40, 4).NumberFormat = "@"
$ws.Cells.Item(40, 4).Value = "124.95"
$ws.Cells.Item(40, 4).ClearFormats()
$ws.Range("E40").Value = "  +3.99%  "
$ws.Cells.Item(41, 4).NumberFormat = "@"
$ws.Cells.Item(41, 4).Value = "22.54"
$ws.Cells.Item(41, 4).ClearFormats()
$ws.Range("E41").Value = "  +2.39%  "
$ws.Cells.Item(42, 4).NumberFormat = "@"
$ws.Cells.Item(42, 4).Value = "0.110"
$ws.Cells.Item(42, 4).ClearFormats()
$ws.Range("E42").Value = "  -0.98%  "
$ws.Range("E43").Value = "  -3.38%  "
$ws.Cells.Item(44, 4).NumberFormat = "@"
$ws.Cells.Item(44, 4).Value = "0.0314"
$ws.Cells.Item(44, 4).ClearFormats()
$ws.Range("E44").Value = "  +1.46%  "
$ws.Cells.Item(45, 4).NumberFormat = "@"
$ws.Cells.Item(45, 4).Value = "2.064.46"
$ws.Cells.Item(45, 4).ClearFormats()
$ws.Range("E45").Value = "  +2.42%  "
$ws.Range("E46").Value = "  -1.35%  "
$ws.Range("E47").Value = "  +6.81%  "
$ws.Cells.Item(48, 4).NumberFormat = "@"
$ws.Cells.Item(48, 4).Value = "2.17"
$ws.Cells.Item(48, 4).ClearFormats()
$ws.Range("E48").Value = "  +2.36%  "
$ws.Cells.Item(49, 4).NumberFormat = "@"
$ws.Cells.Item(49, 4).Value = "8.98"
$ws.Cells.Item(49, 4).ClearFormats()
$ws.Range("E49").Value = "  -0.58%  "
$ws.Cells.Item(50, 4).NumberFormat = "@"
$ws.Cells.Item(50, 4).Value = "58.62"
$ws.Cells.Item(50, 4).ClearFormats()
$ws.Range("E50").Value = "  +2.68%  "
$ws.Cells.Item(51, 4).NumberFormat = "@"
$ws.Cells.Item(51, 4).Value = "5.17"
$ws.Cells.Item(51, 4).ClearFormats()
$ws.Range("E51").Value = "  -2.93%  "
